$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 11 with the new "Chocolate Pickup" / "Recurrssion" entry
$ws.Range("A11").Value = "CN"
$ws.Range("B11").Value = "CN"
$ws.Range("C11").Value = "Chocolate Pickup"
$ws.Range("D11").Value = "Java"
$ws.Range("E11").Value = "Recurrssion"

# Match the cell formatting used by neighbouring wrap-text cells in these columns
# (copy the exact style already used by C21/E2, rather than re-deriving it)
$ws.Range("C21").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null

$ws.Range("E2").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Move the active selection to C11 (matches the saved view state in the diff)
$ws.Range("C11").Select()
